$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Widen column D to fit the new, longer "Name" entry (closest the host lets
# us get to the authored 22.8166-character width)
$ws.Columns.Item(4).ColumnWidth = 22

# Add new issue #6 row
# Note: set F7 (Description) before D7 (Name) so new shared strings are
# appended in the same order as the target workbook (Description text
# first, then the short Name text).
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "DONE"
$ws.Range("F7").Value = "Make buttons on image browser not full width … like the ones on the home page"
$ws.Range("D7").Value = "Tidy image browser buttons"
$ws.Range("E7").Value = "Tidy UI"

# The row wraps onto two lines (same as the other Story/Status rows), so it
# needs the same row height as those
$ws.Rows.Item(7).RowHeight = 29

# Move the selection to the newly added cell, matching the author's last edit point
$ws.Range("D7").Select()
